# Equip.xlsx — "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to be a generic "Property1" tab is being folded into
# the unified DataNode/DataTable/Entity model, so it gets renamed to
# "DataNode". The author's last on-screen selection (in the frozen
# bottom-left pane) is also captured as part of the same save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Re-create the author's final selection: top pane stayed on I1, the
# frozen bottom-left pane's active cell moved to L39.
$ws.Range("I1").Select()
$ws.Range("L39").Select()
